$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: update Correspond Handoff/Handback Datetime values ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-12 00:15:03"
$wsZh.Range("H3").Value = "2016-03-12 00:15:44"
$wsZh.Range("E5").Value = "2016-03-12 00:15:03"
$wsZh.Range("H5").Value = "2016-03-12 00:15:44"

# --- de-de sheet: update Correspond Handoff/Handback Datetime values ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-12 00:15:07"
$wsDe.Range("H3").Value = "2016-03-12 00:15:53"
$wsDe.Range("E5").Value = "2016-03-12 00:15:07"
$wsDe.Range("H5").Value = "2016-03-12 00:15:53"
